$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 25.02674933333333
$ws.Range("H2").Value = 75.080248
$ws.Range("I2").Value = 0.05787790829091637
$ws.Range("J2").Value = 0.05787790829091637
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 98.946724
$ws.Range("N2").Value = 296.840172
$ws.Range("O2").Value = 0.2098009692989996
$ws.Range("P2").Value = 0.2098009692989996
$ws.Range("Q2").Value = 2476.314858902517
$ws.Range("R2").Value = 22286.83373012266
$ws.Range("S2").Value = 0.01214284126043286
$ws.Range("T2").Value = 0.01214284126043286
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 25.02674933333333
$ws.Range("H3").Value = 75.080248
$ws.Range("I3").Value = 0.05787790829091637
$ws.Range("J3").Value = 0.05787790829091637
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 163.0062356666667
$ws.Range("N3").Value = 489.018707
$ws.Range("O3").Value = 0.345629090707923
$ws.Range("P3").Value = 0.3456290907079231
$ws.Range("Q3").Value = 4079.516199799926
$ws.Range("R3").Value = 36715.64579819934
$ws.Range("S3").Value = 0.02000428881466599
$ws.Range("T3").Value = 0.02000428881466599
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 25.02674933333333
$ws.Range("H4").Value = 75.080248
$ws.Range("I4").Value = 0.05787790829091637
$ws.Range("J4").Value = 0.05787790829091637
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 65.39610666666668
$ws.Range("N4").Value = 196.18832
$ws.Range("O4").Value = 0.1386621609326595
$ws.Range("P4").Value = 0.1386621609326595
$ws.Range("Q4").Value = 1636.651968922596
$ws.Range("R4").Value = 14729.86772030336
$ws.Range("S4").Value = 0.008025475833880752
$ws.Range("T4").Value = 0.008025475833880752
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 25.02674933333333
$ws.Range("H5").Value = 75.080248
$ws.Range("I5").Value = 0.05787790829091637
$ws.Range("J5").Value = 0.05787790829091637
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 144.2727966666667
$ws.Range("N5").Value = 432.81839
$ws.Range("O5").Value = 0.3059077790604178
$ws.Range("P5").Value = 0.3059077790604179
$ws.Range("Q5").Value = 3610.679117795635
$ws.Range("R5").Value = 32496.11206016072
$ws.Range("S5").Value = 0.01770530238193677
$ws.Range("T5").Value = 0.01770530238193677
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 138.9376323333333
$ws.Range("H6").Value = 416.812897
$ws.Range("I6").Value = 0.3213129853678317
$ws.Range("J6").Value = 0.3213129853678316
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 98.946724
$ws.Range("N6").Value = 296.840172
$ws.Range("O6").Value = 0.2098009692989996
$ws.Range("P6").Value = 0.2098009692989996
$ws.Range("Q6").Value = 13747.42355969981
$ws.Range("R6").Value = 123726.8120372983
$ws.Range("S6").Value = 0.06741177577852636
$ws.Range("T6").Value = 0.06741177577852636
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 138.9376323333333
$ws.Range("H7").Value = 416.812897
$ws.Range("I7").Value = 0.3213129853678317
$ws.Range("J7").Value = 0.3213129853678316
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.0062356666667
$ws.Range("N7").Value = 489.018707
$ws.Range("O7").Value = 0.345629090707923
$ws.Range("P7").Value = 0.3456290907079231
$ws.Range("Q7").Value = 22647.70043909602
$ws.Range("R7").Value = 203829.3039518642
$ws.Range("S7").Value = 0.1110551149653319
$ws.Range("T7").Value = 0.1110551149653319
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 138.9376323333333
$ws.Range("H8").Value = 416.812897
$ws.Range("I8").Value = 0.3213129853678317
$ws.Range("J8").Value = 0.3213129853678316
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 65.39610666666668
$ws.Range("N8").Value = 196.18832
$ws.Range("O8").Value = 0.1386621609326595
$ws.Range("P8").Value = 0.1386621609326595
$ws.Range("Q8").Value = 9085.980224084786
$ws.Range("R8").Value = 81773.82201676306
$ws.Range("S8").Value = 0.04455395288682754
$ws.Range("T8").Value = 0.04455395288682754
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 138.9376323333333
$ws.Range("H9").Value = 416.812897
$ws.Range("I9").Value = 0.3213129853678317
$ws.Range("J9").Value = 0.3213129853678316
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 144.2727966666667
$ws.Range("N9").Value = 432.81839
$ws.Range("O9").Value = 0.3059077790604178
$ws.Range("P9").Value = 0.3059077790604179
$ws.Range("Q9").Value = 20044.92077897509
$ws.Range("R9").Value = 180404.2870107758
$ws.Range("S9").Value = 0.09829214173714593
$ws.Range("T9").Value = 0.09829214173714593
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 98.29468800000001
$ws.Range("H10").Value = 294.884064
$ws.Range("I10").Value = 0.2273204107243322
$ws.Range("J10").Value = 0.2273204107243321
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 98.946724
$ws.Range("N10").Value = 296.840172
$ws.Range("O10").Value = 0.2098009692989996
$ws.Range("P10").Value = 0.2098009692989996
$ws.Range("Q10").Value = 9725.937364202113
$ws.Range("R10").Value = 87533.43627781901
$ws.Range("S10").Value = 0.04769204251141158
$ws.Range("T10").Value = 0.04769204251141158
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 98.29468800000001
$ws.Range("H11").Value = 294.884064
$ws.Range("I11").Value = 0.2273204107243322
$ws.Range("J11").Value = 0.2273204107243321
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 163.0062356666667
$ws.Range("N11").Value = 489.018707
$ws.Range("O11").Value = 0.345629090707923
$ws.Range("P11").Value = 0.3456290907079231
$ws.Range("Q11").Value = 16022.64707690947
$ws.Range("R11").Value = 144203.8236921853
$ws.Range("S11").Value = 0.07856854685800252
$ws.Range("T11").Value = 0.07856854685800251
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 98.29468800000001
$ws.Range("H12").Value = 294.884064
$ws.Range("I12").Value = 0.2273204107243322
$ws.Range("J12").Value = 0.2273204107243321
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 65.39610666666668
$ws.Range("N12").Value = 196.18832
$ws.Range("O12").Value = 0.1386621609326595
$ws.Range("P12").Value = 0.1386621609326595
$ws.Range("Q12").Value = 6428.089901214722
$ws.Range("R12").Value = 57852.8091109325
$ws.Range("S12").Value = 0.03152073937513559
$ws.Range("T12").Value = 0.03152073937513559
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 98.29468800000001
$ws.Range("H13").Value = 294.884064
$ws.Range("I13").Value = 0.2273204107243322
$ws.Range("J13").Value = 0.2273204107243321
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 144.2727966666667
$ws.Range("N13").Value = 432.81839
$ws.Range("O13").Value = 0.3059077790604178
$ws.Range("P13").Value = 0.3059077790604179
$ws.Range("Q13").Value = 14181.24953523744
$ws.Range("R13").Value = 127631.245817137
$ws.Range("S13").Value = 0.06953908197978244
$ws.Range("T13").Value = 0.06953908197978242
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 170.1468356666666
$ws.Range("H14").Value = 510.440507
$ws.Range("I14").Value = 0.3934886956169198
$ws.Range("J14").Value = 0.3934886956169198
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 98.946724
$ws.Range("N14").Value = 296.840172
$ws.Range("O14").Value = 0.2098009692989996
$ws.Range("P14").Value = 0.2098009692989996
$ws.Range("Q14").Value = 16835.47198818302
$ws.Range("R14").Value = 151519.2478936472
$ws.Range("S14").Value = 0.08255430974862879
$ws.Range("T14").Value = 0.08255430974862879
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 170.1468356666666
$ws.Range("H15").Value = 510.440507
$ws.Range("I15").Value = 0.3934886956169198
$ws.Range("J15").Value = 0.3934886956169198
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 163.0062356666667
$ws.Range("N15").Value = 489.018707
$ws.Range("O15").Value = 0.345629090707923
$ws.Range("P15").Value = 0.3456290907079231
$ws.Range("Q15").Value = 27734.99519261827
$ws.Range("R15").Value = 249614.9567335644
$ws.Range("S15").Value = 0.1360011400699227
$ws.Range("T15").Value = 0.1360011400699227
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 170.1468356666666
$ws.Range("H16").Value = 510.440507
$ws.Range("I16").Value = 0.3934886956169198
$ws.Range("J16").Value = 0.3934886956169198
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 65.39610666666668
$ws.Range("N16").Value = 196.18832
$ws.Range("O16").Value = 0.1386621609326595
$ws.Range("P16").Value = 0.1386621609326595
$ws.Range("Q16").Value = 11126.94061425314
$ws.Range("R16").Value = 100142.4655282782
$ws.Range("S16").Value = 0.05456199283681559
$ws.Range("T16").Value = 0.05456199283681559
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 170.1468356666666
$ws.Range("H17").Value = 510.440507
$ws.Range("I17").Value = 0.3934886956169198
$ws.Range("J17").Value = 0.3934886956169198
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 144.2727966666667
$ws.Range("N17").Value = 432.81839
$ws.Range("O17").Value = 0.3059077790604178
$ws.Range("P17").Value = 0.3059077790604179
$ws.Range("Q17").Value = 24547.55982561375
$ws.Range("R17").Value = 220928.0384305237
$ws.Range("S17").Value = 0.1203712529615527
$ws.Range("T17").Value = 0.1203712529615527
